$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump version and publication date ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B3").Value = "2.2.0-ballot"
$metaSheet.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# --- Include ValueSet #0: pin version on the referenced ValueSet URL ---
$vs0Sheet = $wb.Worksheets.Item("Include ValueSet #0")
$vs0Sheet.Range("A2").Value = "https://hl7.fr/ig/fhir/core/ValueSet/fr-core-vs-encounter-type|2.1.0"

# --- Include ValueSet #2: pin version on the referenced ValueSet URL ---
$vs2Sheet = $wb.Worksheets.Item("Include ValueSet #2")
$vs2Sheet.Range("A2").Value = "https://smt.esante.gouv.fr/fhir/ValueSet/jdv-type-evenement-ssiad-cisis|20250624152100"
